$d = $word.ActiveDocument

# --- 1. Add the new "Week9" progress paragraph right after the "Week7" one ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "Week7:*") {
        $targetIdx = $i
    }
}

if ($targetIdx -gt 0) {
    $week7Para = $d.Paragraphs.Item($targetIdx)
    $week7Para.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.Text = "Week9: Working on segnet model to perform the segmentation of iris."
}

# --- 2. Update the Heading1-6 styles ---
$headingBase = $d.Styles.Item("Heading")

foreach ($name in @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6")) {
    $s = $d.Styles.Item($name)
    $s.BaseStyle = $headingBase
    $s.ParagraphFormat.ReadingOrder = 0
    $s.ParagraphFormat.Alignment = 0
}

# Heading 1 and Heading 2 also pick up the darker explicit text color.
$d.Styles.Item("Heading 1").Font.Color = 655360
$d.Styles.Item("Heading 2").Font.Color = 655360

Write-Output "ok"
